# Add the 29/03/2020 date column (Q) to the DIY infected-cases tracker.
# Column P currently holds 28/03/2020 as the last (specially-styled) date
# column; R is a blank spacer. We shift that "last column" styling from P
# onto the new Q column, restyle P back to a normal data column, give the
# new Q header cell the filled title-row look used across row 1, and fill
# the new date's case counts with 0 for every Kecamatan row (3-81), mirroring
# the existing column P values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 header: new date label in Q2, carrying P2's current (accent) format ---
$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)
$ws.Range("Q2").Value2 = "29/03/2020"

# P2 reverts to the plain header format used by the other date columns (e.g. O2)
$ws.Range("O2").Copy()
$ws.Range("P2").PasteSpecial(-4122)

# --- Row 1 title band: Q1 picks up the filled title style used across A1:P1 ---
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)

# --- Data rows 3-81: new Q column, seeded with 0 like the rest of the table ---
$ws.Range("P3:P81").Copy()
$ws.Range("Q3:Q81").PasteSpecial(-4122)
$ws.Range("Q3:Q81").Value2 = 0.0

$excel.CutCopyMode = $false
